# Update "想去人数" (F column) counts by +1 for specific events
# on both the "展览" sheet and the "全部类型" sheet, mirroring the
# site regeneration output captured in the commit.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAll = $wb.Worksheets.Item("全部类型")

# Row numbers on "展览" sheet with their new F-column values
$exhibitionUpdates = @{
    2  = 3128
    10 = 15578
    14 = 6143
    24 = 10
    26 = 208
    28 = 23
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Row numbers on "全部类型" sheet with their new F-column values
$allTypesUpdates = @{
    3  = 3128
    11 = 15578
    15 = 6143
    25 = 10
    27 = 208
    29 = 23
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
